$p = $ppt.ActivePresentation
$p2 = $ppt.Presentations.Open("before.pptx")
Write-Host "p2 designs count:" $p2.Designs.Count
$d2 = $p2.Designs.Item(1)
Write-Host "d2 name:" $d2.Name
try {
  $newD = $p.Designs.Clone($d2, 2)
  Write-Host "Cloned:" $newD.Name $newD.Index
} catch {
  Write-Host "ERROR:" $_.Exception.Message
}
Write-Host "p designs count:" $p.Designs.Count
